# Updated dcp2T1 with Amnon's notebook
# Applies updates to the "Tier 1_obs" sheet for data rows 6-34:
#   - assay_ontology_term_id (Z)  : EFO:0009900        -> EFO:0009900||EFO:0010714
#   - reference_genome (AF)       : GRCh37              -> GRCh37||GRCh38
#   - cell_enrichment (V)         : na                  -> DAPI-||nan   (all data rows except 6, 11, 16, 22)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

# Rows that keep their existing cell_enrichment value ("na") unchanged.
$skipV = @(6, 11, 16, 22)

for ($row = 6; $row -le 34; $row++) {
    $ws.Range("Z$row").Value = "EFO:0009900||EFO:0010714"
    $ws.Range("AF$row").Value = "GRCh37||GRCh38"

    if ($skipV -notcontains $row) {
        $ws.Range("V$row").Value = "DAPI-||nan"
    }
}
